# Find a branch test case removal from PROD(TC19) and ECTEST(TC21)
# Remove the "TC21_Verify_ Find_a_Branch_Guestuser" row (row 22) from the
# ECTEST MasterExecutor_Sanity sheet; everything below shifts up.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(22).Delete()

# Re-apply the AutoFilter over the now-smaller data range and refresh the
# hidden _FilterDatabase defined name so it matches the new extents.
$ws.Range("A1:F29").AutoFilter()

foreach ($n in $wb.Names) {
    if ($n.Name -eq "MasterExecutor!_FilterDatabase") {
        $n.RefersTo = "=MasterExecutor!`$A`$1:`$F`$29"
    }
}

# Reflect the row-22 selection left behind by the deletion.
$ws.Rows.Item(22).Select()
